$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Simple single-value cell replacements
Replace-Exact "65.87" "0M"
Replace-Exact "214.66" "0M"
Replace-Exact "628" "0M"
Replace-Exact "10438" "12135"
Replace-Exact "0.07408" "0.07452"
Replace-Exact "0.02071" "0.02199"
Replace-Exact "0.01167" "0.01232"
Replace-Exact "0.01676" "0.01683"
Replace-Exact "0.02965" "0.02983"
Replace-Exact "166.72086" "214.65973"

# Collapse the tab-separated multi-value cells (rows 44-46) down to a
# single value each, dropping the extra tab-delimited figures.
$tbl = $d.Tables.Item(1)
$tbl.Cell(44, 1).Range.Text = "65.87"
$tbl.Cell(45, 1).Range.Text = "214.66"
$tbl.Cell(46, 1).Range.Text = "628"
